# Update the error-structure triangle (Datos_sal/estructura_errores_mod2)
# Every populated cell in the triangular data range (A2:L25) is replaced
# with its refreshed value, row by row, leaving headers (row 1) and the
# empty cells outside the triangle untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @('A','B','C','D','E','F','G','H','I','J','K','L')

$rowValues = @{}
$rowValues[2] = @(0.000189094899743836, 0.00229290042824465, 0.0047509171634274, 0.0065001039743809, 0.00715736812945828, 0.00780661786451392, 0.00834796155302654, 0.00946250465293428, 0.0112913479336969, 0.0128646005310842, 0.0142102995597854, 0.0151900860648983)
$rowValues[3] = @(0.00288275075713695, 0.00541184884135097, 0.00705360428833929, 0.00751858201974164, 0.00805143918831408, 0.00850960370263281, 0.00961634013915217, 0.0114978964827267, 0.0130939925495329, 0.0144441589832385, 0.0154078485869394, 0.0160733214295677)
$rowValues[4] = @(0.00172247832346528, 0.00209005154835023, 0.0017491146548003, 0.00167421650414747, 0.00175996162774325, 0.00306765231712405, 0.00546578531339974, 0.00721290458053446, 0.0086086381293833, 0.00952906207874279, 0.0101128332239047, 0.0103483454229489)
$rowValues[5] = @(0.000793843575943143, 0.00257395536880098, 0.00285621251453618, 0.00287050386435736, 0.00264627842777143, 0.00390525477514005, 0.00514885722627092, 0.00625087050944847, 0.00694905110375476, 0.00735172855245885, 0.00743434284528481, 0.00720658971617156)
$rowValues[6] = @(0.00208222381328493, 0.00177879710173363, 0.00151020792065469, 0.00227339149017842, 0.00482646842446357, 0.0065603589274265, 0.00790948044871535, 0.0087325945869329, 0.00920157862453353, 0.00931018914517136, 0.00905495205507859, 0.00873944221221079)
$rowValues[7] = @(0.00244238011228237, 0.00353318343485633, 0.00617525001072711, 0.00951957595559136, 0.0116637632066486, 0.0132099055200919, 0.0141206310863958, 0.0146192597049183, 0.0147229781046165, 0.0144074053289671, 0.0140066054453225, 0.01345116007408)
$rowValues[8] = @(0.000151465283849461, 0.0025015091130499, 0.00561083144581903, 0.00728626696195184, 0.00850976591477936, 0.00917284549851884, 0.00949232425581297, 0.00945850044666312, 0.00907066274101104, 0.0086549378136001, 0.00828671601926887, 0.00800684040680302)
$rowValues[9] = @(0.0038145431891321, 0.00718610267833207, 0.00876378455386702, 0.00988601525963332, 0.0104324394339922, 0.0106429711172854, 0.0105034062973565, 0.0100069526014997, 0.00949920117944713, 0.00903673469135935, 0.00867524890977546, 0.00834654285967573)
$rowValues[10] = @(0.00238973081235283, 0.00221310712058812, 0.00221012317179192, 0.0020553291658927, 0.00187561714530274, 0.00175038633322698, 0.00242962838359355, 0.00324764682658814, 0.00498922715509503, 0.00625816946318555, 0.00704447744959054, 0.00713954668684147)
$rowValues[11] = @(0.00239186719590112, 0.00305884993311912, 0.00385554645859847, 0.00449957192789146, 0.00530016577101576, 0.00672297185472445, 0.00792748969956822, 0.00978637363919558, 0.0112055326084795, 0.0121350332883279, 0.012336159469511, 0.0120312594724519)
$rowValues[12] = @(0.000827047829735861, 0.000787976573667292, 0.00072328259744124, 0.000750728850677807, 0.00216088289090913, 0.00325311464401535, 0.00532802186703687, 0.00673756589088804, 0.00757494079782967, 0.00763073645995737, 0.00729749377662428, 0.00699531571125448)
$rowValues[13] = @(0.000782994920796476, 0.00116087653911119, 0.00202970681598535, 0.00390184248666624, 0.00522500241116569, 0.00741973460088451, 0.00891669080687451, 0.0097967400806848, 0.00984573430317491, 0.00942432632896037, 0.00904022347033389, 0.00866667296111902)
$rowValues[14] = @(0.00000181763959616887, 0.000855522361654251, 0.00286618738997368, 0.00406989845807709, 0.00637054570000637, 0.00783763209918939, 0.00864496290146159, 0.00858912839495094, 0.00812418551992946, 0.00771815930735211, 0.00736027672024591, 0.00705185170120391)
$rowValues[15] = @(0.00121331060299656, 0.00351404821421334, 0.00470373852500934, 0.0071268663411251, 0.00859030230786514, 0.00934239474918559, 0.0091870117778295, 0.00862166702991275, 0.00814013682410728, 0.00772374658906264, 0.00736949600024266)
$rowValues[16] = @(0.0025793236691678, 0.00307889479728456, 0.005442442654035, 0.00665825965385098, 0.00718635656526592, 0.00685181810803247, 0.00636167476813394, 0.00598622708963998, 0.00576540598653844, 0.00551100206410359)
$rowValues[17] = @(0.00126029812888717, 0.0017274990441961, 0.00194508470716094, 0.00186823980426013, 0.00200471977723274, 0.00388286129979889, 0.0048747021975488, 0.00581470811675162, 0.00619816514169062)
$rowValues[18] = @(0.00442246525274292, 0.00492683909227809, 0.00497752347247237, 0.00434725549922881, 0.00444452641028206, 0.00458718181546982, 0.0049754719252769, 0.00503438436624933)
$rowValues[19] = @(0.00278346689398656, 0.00447020468857089, 0.0073484481091654, 0.0108121156428032, 0.012535329112032, 0.0138425190707409, 0.0143734315933842)
$rowValues[20] = @(0.000541070704052515, 0.00305970230630726, 0.00620122564659881, 0.00728920868554469, 0.00820548354668938, 0.00843473047165011)
$rowValues[21] = @(0.00329710697071928, 0.00645982086496855, 0.00715545168055173, 0.00783799240996244, 0.0078640584769839)
$rowValues[22] = @(0.00246216774691632, 0.00176632211134729, 0.00148474863522086, 0.00145760430693607)
$rowValues[23] = @(0.00412951484759283, 0.00479929389008437, 0.00612707790199198)
$rowValues[24] = @(0.00222686008955808, 0.00204569298691562)
$rowValues[25] = @(0.00225059521058402)

foreach ($r in ($rowValues.Keys | Sort-Object)) {
    $vals = $rowValues[$r]
    $n = $vals.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $lastCol = $colLetters[$n-1]
    $ws.Range("A$($r):$($lastCol)$($r)").Value2 = $arr
}

